$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Salary") transaction: amount and date updated
$ws.Range("B2").Value = 120000
$ws.Range("C2").Value = 46094.22928240741

# Row 3 ("Salary") transaction: amount and date updated
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 46024.22928240741

# Row 4 transaction removed entirely (shifts remaining rows up)
$ws.Rows.Item(4).Delete()
